$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.438.03"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.885.00"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.61"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4843"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2895"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06601"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Value = "1.890.83"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.87"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07411"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.179"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.46"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6619"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "30.432.47"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.54"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007750"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "2.139.79"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.384"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "218.28"
$ws.Range("E23").Value = "  +14.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.180"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.373"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.23"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.62"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.936"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.445"
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.338"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09181"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.038"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05073"
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7544"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.710"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01888"
$ws.Range("E37").Value = "  +3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.647"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9205"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.095"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.990"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.75"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4343"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.668"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1333"
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.591"
$ws.Range("E47").Value = "  +11.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.29"
$ws.Range("E48").Value = "  -12.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.867"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.48"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05710"
$ws.Range("E51").Value = "  -2.44%  "
